# Upgrade left table until 2023 (Javakheti diff context): add a new "2023"
# column (K) to the Sachkhere average monthly remuneration table, mirroring
# the formatting already used for the 2022 column (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number formatting/styles from the 2022 column (J, rows 3-6,
# which hold the year header and the three data rows) onto the new 2023
# column (K) before filling in values, so the new cells pick up the same
# cell styles (s="5", s="11", s="12", s="13") as column J.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Fill in the new 2023 data.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1198.7
$ws.Range("K5").Value = 914.9
$ws.Range("K6").Value = 1420.5
